# Append the two new "kilométrage" readings that were recorded after the
# last existing row (row 76, 2020-04-09 / 4713 km).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 77: 2020-04-20 -> 4726 km
$ws.Cells.Item(77, 1).Value = 43941
$ws.Cells.Item(77, 2).Value = 4726

# Row 78: 2020-04-25 -> 4735 km
$ws.Cells.Item(78, 1).Value = 43946
$ws.Cells.Item(78, 2).Value = 4735

# Mirror the author's UI state: scrolled down so row 50 is at the top, with
# the next empty cell (B79) selected, ready for the following entry.
$win = $excel.ActiveWindow
$win.ScrollRow = 50
$win.ScrollColumn = 1
$ws.Range("B79").Select()
